$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1. "positions" sheet tweaks
#    - clear the start_date for position 3 (COO / Charlie)
#    - set an end_date for position 4 (Sales Rep / Dave), using the m/d/yy
#      date style already used for end_date elsewhere (e.g. related!J2)
# ---------------------------------------------------------------------------
$positions = $wb.Worksheets.Item("positions")
$related = $wb.Worksheets.Item("related")

$positions.Range("I4").ClearContents()

$related.Range("J2").Copy()
$positions.Range("J5").PasteSpecial($xlPasteFormats)
$positions.Range("J5").Value = 44727

# ---------------------------------------------------------------------------
# 2. Split the "inflation" sheet in two:
#    - duplicate it so the formulas/content are preserved on a fresh sheet
#    - rename the original to "onetime" and replace its contents with the
#      new one-time personnel expense data
#    - rename the duplicate back to "inflation"
# ---------------------------------------------------------------------------
$inflation = $wb.Worksheets.Item("inflation")
$inflation.Copy($null, $inflation)
$inflation.Name = "onetime"
$wb.Worksheets.Item("inflation (2)").Name = "inflation"

$onetime = $wb.Worksheets.Item("onetime")
$onetime.Cells.Clear()

# ---- header row ----
$onetime.Cells.Item(1, 1).Value = "position_id"
$onetime.Cells.Item(1, 2).Value = "position_title"
$onetime.Cells.Item(1, 3).Value = "department"
$onetime.Cells.Item(1, 4).Value = "employee_id"
$onetime.Cells.Item(1, 5).Value = "employee_name"
$onetime.Cells.Item(1, 6).Value = "item"
$onetime.Cells.Item(1, 7).Value = "expense_type"
$onetime.Cells.Item(1, 8).Value = "expense_amount"
$onetime.Cells.Item(1, 9).Value = "expense_date"

# ---- row 2: signing bonus for Bob (CFO, position 2) ----
$onetime.Cells.Item(2, 1).Value = 2
$onetime.Cells.Item(2, 2).Value = "CFO"
$onetime.Cells.Item(2, 3).Value = "d2"
$onetime.Cells.Item(2, 4).Value = "e10002"
$onetime.Cells.Item(2, 5).Value = "Bob"
$onetime.Cells.Item(2, 6).Value = "signing bonud"
$onetime.Cells.Item(2, 7).Value = "salary"
$onetime.Cells.Item(2, 8).Value = 10000
$onetime.Cells.Item(2, 9).Value = 44607

# ---- row 3: severance for Dave (Sales Rep, position 4) ----
$onetime.Cells.Item(3, 1).Value = 4
$onetime.Cells.Item(3, 2).Value = "Sales Rep"
$onetime.Cells.Item(3, 3).Value = "d4"
$onetime.Cells.Item(3, 4).Value = "e10004"
$onetime.Cells.Item(3, 5).Value = "Dave"
$onetime.Cells.Item(3, 6).Value = "severance"
$onetime.Cells.Item(3, 7).Value = "salary"
$onetime.Cells.Item(3, 8).Value = 5000
$onetime.Cells.Item(3, 9).Value = 44727

# ---- formatting: reuse the styles already used elsewhere in the workbook ----

# header row -> bold header style (positions!A1:I1)
$positions.Range("A1:I1").Copy()
$onetime.Range("A1:I1").PasteSpecial($xlPasteFormats)

# row2/row3 text & id columns -> plain style (positions!A2:E2)
$positions.Range("A2:E2").Copy()
$onetime.Range("A2:E3").PasteSpecial($xlPasteFormats)

# "item"/"expense_type" + currency "expense_amount" columns -> currency style (positions!F2)
$positions.Range("F2").Copy()
$onetime.Range("F2:F3").PasteSpecial($xlPasteFormats)
$onetime.Range("H2:H3").PasteSpecial($xlPasteFormats)

# "expense_type" column -> plain style (positions!A2)
$positions.Range("A2").Copy()
$onetime.Range("G2:G3").PasteSpecial($xlPasteFormats)

# "expense_date" column -> m/d/yyyy date style (positions!I3)
$positions.Range("I3").Copy()
$onetime.Range("I2:I3").PasteSpecial($xlPasteFormats)

# trailing, still-empty cells that keep formatting in the source workbook
$positions.Range("I2").Copy()
$onetime.Range("J2:J3").PasteSpecial($xlPasteFormats)

$related.Range("J2").Copy()
$onetime.Range("K3").PasteSpecial($xlPasteFormats)

# restore the values (PasteSpecial(xlPasteFormats) should not touch values, but
# make sure by re-asserting them after all formatting is applied)
$onetime.Cells.Item(2, 9).Value = 44607
$onetime.Cells.Item(3, 9).Value = 44727

$onetime.Range("A1").Select()
